$d = $word.ActiveDocument

$replacements = @(
    @{old="2024-04-11 Thursday"; new="2024-04-12 Friday"},
    @{old="976×2="; new="575×9="},
    @{old="275×6="; new="928×9="},
    @{old="553×8="; new="623×6="},
    @{old="832×8="; new="503×5="},
    @{old="899×2="; new="969×6="},
    @{old="737×2="; new="674×6="},
    @{old="769×4="; new="158×6="},
    @{old="548×3="; new="677×4="},
    @{old="339×8="; new="391×4="},
    @{old="438×2="; new="487×4="},
    @{old="950×8="; new="815×6="},
    @{old="278×9="; new="474×8="},
    @{old="500×4="; new="690×9="},
    @{old="775×9="; new="719×9="},
    @{old="448×7="; new="620×6="},
    @{old="429×5="; new="327×6="},
    @{old="727×2="; new="431×3="},
    @{old="231×6="; new="332×9="},
    @{old="660×9="; new="532×4="},
    @{old="895×3="; new="915×7="},
    @{old="814×4="; new="346×4="},
    @{old="618×5="; new="704×2="},
    @{old="490×8="; new="434×9="},
    @{old="730×2="; new="291×2="},
    @{old="196×8="; new="557×2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
